# QA Excel Compiler update: add STATUS_{Username} columns, reorder
# COMMENT_{Username} / COMMENT columns, and refresh sample STATUS +
# COMMENT sample data on Sheet1 / Sheet2, plus recompute the Alice
# completion percentages on the STATUS sheet.

$wb = $excel.ActiveWorkbook

function Set-QuestSheetHeader($ws) {
    $ws.Cells.Item(1, 1).Value = "Original"
    $ws.Cells.Item(1, 2).Value = "ENG"
    $ws.Cells.Item(1, 3).Value = "StringKey"
    $ws.Cells.Item(1, 4).Value = "Command"
    $ws.Cells.Item(1, 5).Value = "STATUS"
    $ws.Cells.Item(1, 6).Value = "STATUS_John"
    $ws.Cells.Item(1, 7).Value = "STATUS_Bob"
    $ws.Cells.Item(1, 8).Value = "STATUS_Alice"
    $ws.Cells.Item(1, 9).Value = "COMMENT_John"
    $ws.Cells.Item(1, 10).Value = "COMMENT_Bob"
    $ws.Cells.Item(1, 11).Value = "COMMENT_Alice"
    $ws.Cells.Item(1, 12).Value = "COMMENT"
    $ws.Cells.Item(1, 13).Value = "SCREENSHOT"
}

# ---------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

Set-QuestSheetHeader $ws1

# Row 2
$ws1.Cells.Item(2, 6).Value = "NO ISSUE"
$ws1.Cells.Item(2, 7).Value = ""
$ws1.Cells.Item(2, 8).Value = "NO ISSUE"
$ws1.Cells.Item(2, 9).Value = '"Translation looks good" (date: 251230 1219)'

# Row 3
$ws1.Cells.Item(3, 6).Value = "ISSUE"
$ws1.Cells.Item(3, 7).Value = "ISSUE"
$ws1.Cells.Item(3, 8).Value = "NO ISSUE"
$ws1.Cells.Item(3, 9).Value = '"Typo: should be singular" (date: 251230 1219)'
$ws1.Cells.Item(3, 10).Value = '"Agree - typo" (date: 251230 1219)'
$ws1.Cells.Item(3, 11).Value = "`"I think it's fine`" (date: 251230 1219)"

# Row 4
$ws1.Cells.Item(4, 8).Value = "ISSUE"
$ws1.Cells.Item(4, 9).Value = ""
$ws1.Cells.Item(4, 11).Value = '"Missing article" (date: 251230 1219)'

# Row 5
$ws1.Cells.Item(5, 6).Value = "NO ISSUE"
$ws1.Cells.Item(5, 7).Value = "NO ISSUE"
$ws1.Cells.Item(5, 8).Value = ""
$ws1.Cells.Item(5, 9).Value = '"Perfect" (date: 251230 1219)'
$ws1.Cells.Item(5, 10).Value = '"Good" (date: 251230 1219)'

# Row 6
$ws1.Cells.Item(6, 7).Value = "ISSUE"
$ws1.Cells.Item(6, 8).Value = ""
$ws1.Cells.Item(6, 10).Value = '"Forest should be woods" (date: 251230 1219)'

# ---------------------------------------------------------------
# Sheet2
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

Set-QuestSheetHeader $ws2

# Row 2
$ws2.Cells.Item(2, 6).Value = "NO ISSUE"
$ws2.Cells.Item(2, 7).Value = "NO ISSUE"
$ws2.Cells.Item(2, 8).Value = ""
$ws2.Cells.Item(2, 9).Value = '"Verified" (date: 251230 1219)'
$ws2.Cells.Item(2, 10).Value = '"All good" (date: 251230 1219)'

# Row 3
$ws2.Cells.Item(3, 7).Value = "NO ISSUE"
$ws2.Cells.Item(3, 8).Value = "NO ISSUE"
$ws2.Cells.Item(3, 9).Value = ""
$ws2.Cells.Item(3, 10).Value = '"Verified" (date: 251230 1219)'
$ws2.Cells.Item(3, 11).Value = '"Shop verified" (date: 251230 1219)'

# Row 4
$ws2.Cells.Item(4, 6).Value = "BLOCKED"
$ws2.Cells.Item(4, 7).Value = "NO ISSUE"
$ws2.Cells.Item(4, 8).Value = ""
$ws2.Cells.Item(4, 9).Value = '"Need to check context" (date: 251230 1219)'
$ws2.Cells.Item(4, 10).Value = '"Correct" (date: 251230 1219)'

# ---------------------------------------------------------------
# STATUS sheet - recompute Alice's completion percentages
# (leading apostrophe forces text storage, same as the other
# percentage cells on this sheet, instead of being parsed as a
# numeric percentage)
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("STATUS")
$wsStatus.Cells.Item(2, 2).Value = "'60.0%"
$wsStatus.Cells.Item(2, 4).Value = "'46.6%"
